$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 and G3 changed from 1243 to 1245; the dependent formulas in I2:I3
# (=G+H) and M2:M3 (=I*L/1000) recalculate automatically.
$ws.Range("G2").Value = 1245
$ws.Range("G3").Value = 1245

# Move the active cell / selection to I10.
$ws.Range("I10").Select()
